$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

$titleShape = $s1.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange
$tr.Font.Bold = $true
$tr.Font.Color.RGB = 255
